$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) / Volume(1h) (E) for rows 2-48 ---
# D (Price) values are forced to TEXT via NumberFormat="@" so that
# numeric-looking strings (e.g. "1.233") are not auto-converted to
# floating point numbers by Excel -- matching the source data which
# stores these as plain text (inlineStr) cells. The format/style is
# reset back to the default afterwards so no stray style is left on
# the cell.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = "30.379.46"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value2 = "  +2.26%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = "2.093.73"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value2 = "  +0.06%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value2 = "1.002"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value2 = "  -0.92%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "342.71"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  -0.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "1.001"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value2 = "  -0.83%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "0.5229"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value2 = "  +1.43%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = "0.4421"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value2 = "  +1.17%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = "54.60"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value2 = "  +3.01%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "0.09324"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  +1.05%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = "1.168"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value2 = "  +0.60%  "
$ws.Range("E12").Value2 = "  +0.26%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = "8.578"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value2 = "  +3.71%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "6.902"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value2 = "  +2.56%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "2.064.52"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value2 = "  -1.63%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = "101.37"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value2 = "  +2.24%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = "0.00001158"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value2 = "  +0.95%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = "1.003"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value2 = "  -0.80%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "21.13"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value2 = "  +2.01%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "0.06667"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value2 = "  +0.24%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "6.331"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value2 = "  +2.57%  "
$ws.Range("E22").Value2 = "  -0.79%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = "30.393.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value2 = "  +2.22%  "
$ws.Range("E24").Value2 = "  +0.29%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "2.305"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  -0.64%  "
$ws.Range("E26").Value2 = "  -0.13%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "163.03"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  +0.94%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = "2.506"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value2 = "  -0.06%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "133.17"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value2 = "  +0.25%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = "1.137"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value2 = "  +0.87%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "1.658"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value2 = "  +0.98%  "
$ws.Range("E32").Value2 = "  -0.41%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = "6.821"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value2 = "  +9.65%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "6.256"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value2 = "  +1.84%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = "3.854"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value2 = "  -2.13%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "10.13"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value2 = "  -0.43%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value2 = "0.02635"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value2 = "  +3.13%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = "0.06835"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value2 = "  +2.35%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = "0.6981"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value2 = "  +1.95%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "12.58"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value2 = "  +1.31%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "1.339"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  -1.00%  "
$ws.Range("E42").Value2 = "  -0.26%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "0.6806"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value2 = "  +2.34%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "14.34"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value2 = "  +0.10%  "
$ws.Range("E45").Value2 = "  +1.65%  "
$ws.Range("E46").Value2 = "  -0.79%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "1.374"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value2 = "  +18.22%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "3.633"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value2 = "  +0.37%  "

# --- Row 49/50 swapped: BabyDogeCoin <-> ThetaToken, each with new Price/Volume ---
$ws.Range("B49").Value2 = "ThetaToken"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "1.233"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value2 = "  +10.74%  "

$ws.Range("B50").Value2 = "BabyDogeCoin"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "0.00000000341"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -5.10%  "

# --- Row 51 (EOS) price/volume update ---
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = "1.215"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value2 = "  -0.15%  "
